# Update "想去人数" (number of people interested) figures in the F column
# for rows 2-7 on both the "展览" and "全部类型" worksheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    2 = 25
    3 = 1835
    4 = 557
    5 = 1185
    6 = 6142
    7 = 146
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
